$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D contain text-formatted numeric-looking strings (e.g. "488.22", "68.182.89").
# Force text number format before assignment so Excel does not coerce them to floats,
# then restore the default "Normal" style so the cell keeps no explicit style (as in the source).
$dCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price (column D) values
$ws.Range('D2').Value = '68.182.89'
$ws.Range('D3').Value = '3.918.45'
$ws.Range('D5').Value = '488.22'
$ws.Range('D6').Value = '146.50'
$ws.Range('D7').Value = '0.626'
$ws.Range('D9').Value = '0.730'
$ws.Range('D10').Value = '0.172'
$ws.Range('D11').Value = '0.0000356'
$ws.Range('D12').Value = '42.71'
$ws.Range('D13').Value = '10.63'
$ws.Range('D14').Value = '4.543.67'
$ws.Range('D15').Value = '14.77'
$ws.Range('D16').Value = '3.920.92'
$ws.Range('D18').Value = '20.00'
$ws.Range('D20').Value = '68.317.59'
$ws.Range('D21').Value = '443.19'
$ws.Range('D22').Value = '14.73'
$ws.Range('D23').Value = '3.37'
$ws.Range('D24').Value = '88.40'
$ws.Range('D25').Value = '11.70'
$ws.Range('D26').Value = '11.06'
$ws.Range('D27').Value = '3.64'
$ws.Range('D28').Value = '38.84'
$ws.Range('D30').Value = '709.97'
$ws.Range('D31').Value = '13.48'
$ws.Range('D32').Value = '0.130'
$ws.Range('D34').Value = '0.0₃0904'
$ws.Range('D35').Value = '41.36'
$ws.Range('D36').Value = '59.04'
$ws.Range('D37').Value = '5.80'
$ws.Range('D40').Value = '0.380'
$ws.Range('D41').Value = '0.0479'
$ws.Range('D42').Value = '2.85'
$ws.Range('D43').Value = '3.11'
$ws.Range('D47').Value = '3.40'
$ws.Range('D48').Value = '2.14'
$ws.Range('D49').Value = '0.0⁦0343'
$ws.Range('D50').Value = '3.15'
$ws.Range('D51').Value = '144.89'

# Restore default style on the D cells we touched (removes the temporary text format)
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Coin name / link / volume columns (B, C, E) are plain text already
$ws.Range('E2').Value = '  +1.95%  '
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('E5').Value = '  +4.17%  '
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('E10').Value = '  +5.38%  '
$ws.Range('E11').Value = '  +6.37%  '
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('E13').Value = '  +3.09%  '
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('E15').Value = '  -1.61%  '
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('E19').Value = '  -1.73%  '
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('E21').Value = '  +3.36%  '
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('E25').Value = '  +16.77%  '
$ws.Range('E26').Value = '  +15.77%  '
$ws.Range('E27').Value = '  +2.06%  '
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('E29').Value = '  +2.15%  '
$ws.Range('E30').Value = '  -2.56%  '
$ws.Range('E31').Value = '  -1.43%  '
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('E33').Value = '  +2.54%  '
$ws.Range('E34').Value = '  +15.95%  '
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('E36').Value = '  +2.75%  '
$ws.Range('E37').Value = '  +7.78%  '
$ws.Range('E38').Value = '  -4.29%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('E40').Value = '  +13.27%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('E42').Value = '  +14.12%  '
$ws.Range('E43').Value = '  +0.50%  '
$ws.Range('E44').Value = '  +5.57%  '
$ws.Range('E45').Value = '  +1.47%  '
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('E47').Value = '  +0.63%  '
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E49').Value = '  +45.01%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('E50').Value = '  +1.08%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E51').Value = '  +0.94%  '
